# Regenerate the G column ("K" = strikeouts) values for the save-situation
# data rows, replacing the previous "Strike#" derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G.
# Rows 20, 25, and 67 are intentionally omitted because their values are
# unchanged (already 0) per the source diff.
$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 0
    17 = 1
    18 = 0
    19 = 3
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    26 = 2
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 2
    33 = 1
    34 = 0
    35 = 0
    36 = 1
    37 = 4
    38 = 0
    39 = 1
    40 = 2
    41 = 2
    42 = 2
    43 = 0
    44 = 4
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 0
    50 = 2
    51 = 0
    52 = 0
    53 = 2
    54 = 2
    55 = 2
    56 = 2
    57 = 1
    58 = 0
    59 = 0
    60 = 1
    61 = 0
    62 = 0
    63 = 1
    64 = 0
    65 = 0
    66 = 1
    68 = 0
    69 = 1
    70 = 2
    71 = 1
    72 = 0
    73 = 3
    74 = 1
    75 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
